$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Deep Neural Networks review" topic is renamed to "Deep Q-Learning",
# and the old "Neural Network Policies" / "Deep Q-Learning" / "Double Deep
# Q-Learning" sequence collapses into a two-item sequence:
# "Deep Q-Learning studio" and "Double Deep Q-Learning project".
# This shifts every following row's content up by one slot.

$ws.Cells.Item(14, 2).Value = "Deep Q-Learning"
$ws.Cells.Item(15, 2).Value = "Deep Q-Learning"

$ws.Cells.Item(19, 2).Value = "Deep Q-Learning studio"
$ws.Cells.Item(20, 2).Value = "Deep Q-Learning studio"

$ws.Cells.Item(21, 2).Value = "Double Deep Q-Learning project"
$ws.Cells.Item(22, 2).Value = "Double Deep Q-Learning project"
$ws.Cells.Item(23, 2).Value = "Double Deep Q-Learning project"
$ws.Cells.Item(24, 2).Value = "Double Deep Q-Learning project"

$ws.Cells.Item(26, 2).Value = "Policy Optimization Algorithms (PPO)"

$ws.Cells.Item(28, 2).Value = "Implementation of RL using TF-Agents"

# Row height adjustments that accompany the content re-wrap.
$rowsToResize = @(2, 3, 4, 5, 6, 15, 16, 17, 18, 19, 20, 21, 22, 23)
foreach ($r in $rowsToResize) {
    $ws.Rows.Item($r).RowHeight = 20.25
}
